$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.190.67"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.109.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.109.87"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.00%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.158"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.24%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.461"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000228"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.19%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.597.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.146.26"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.100.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "492.97"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.35%  "

$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.45"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.94%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.75"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "26.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.93"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.29%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "57.57"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.46"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "498.89"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.62%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.11"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.24%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.231.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0403"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.59%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0806"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.09%  "

$ws.Range("E41").Value = "  -1.67%  "

$ws.Range("E42").Value = "  +2.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.99%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.261"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.94%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.05%  "

$ws.Range("B47").Value = "PEPE"
$ws.Range("C47").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₃0545"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.68%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.72%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "121.81"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.71%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.68%  "
